# Apply the "Add files via upload" edit to global_results_mnist.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet: remove trailing space
$ws.Name = "Global metrics"

# Row 1: merged title cell
$ws.Range("A1").Value = "Global metrics -> Dataset: mnist"

# Row 2: header row (algorithm names)
$ws.Range("A2").Value = "Test task"
$ws.Range("B2").Value = "Fine-tuning"
$ws.Range("C2").Value = "Joint datasets"
$ws.Range("D2").Value = "Rehearsal 0.1"
$ws.Range("E2").Value = "Rehearsal 0.3"
$ws.Range("F2").Value = "Rehearsal 0.5"
$ws.Range("G2").Value = "EWC"
$ws.Range("H2").Value = "LwF"

# Row 3: Test accuracy on task 1 after task 1
$ws.Range("A3").Value = "Test accuracy on task 1 after task 1"
$ws.Range("B3").Value = 98.2
$ws.Range("C3").Value = 91.5
$ws.Range("D3").Value = 98.58
$ws.Range("E3").Value = 97.90000000000001
$ws.Range("F3").Value = 98.45999999999999
$ws.Range("G3").Value = 98.7
$ws.Range("H3").Value = 98.12

# Row 4: Test accuracy on task 2 after task 1
$ws.Range("A4").Value = "Test accuracy on task 2 after task 1"
$ws.Range("B4").Value = 8.06
$ws.Range("C4").Value = 91.5
$ws.Range("D4").Value = 8.800000000000001
$ws.Range("E4").Value = 6.05
$ws.Range("F4").Value = 6.76
$ws.Range("G4").Value = 8.33
$ws.Range("H4").Value = 4.82

# Row 5: Test average accuracy after task 1
$ws.Range("A5").Value = "Test average accuracy after task 1"
$ws.Range("B5").Value = 53.13
$ws.Range("C5").Value = 91.5
$ws.Range("D5").Value = 53.69
$ws.Range("E5").Value = 51.975
$ws.Range("F5").Value = 52.61
$ws.Range("G5").Value = 53.515
$ws.Range("H5").Value = 51.47

# Row 6: Test accuracy on task 1 after task 2
$ws.Range("A6").Value = "Test accuracy on task 1 after task 2"
$ws.Range("B6").Value = 10.25
$ws.Range("C6").Value = 91.5
$ws.Range("D6").Value = 96.39
$ws.Range("E6").Value = 97.78
$ws.Range("F6").Value = 98.11
$ws.Range("G6").Value = 55.28
$ws.Range("H6").Value = 44.48

# Row 7: Test accuracy on task 2 after task 2
$ws.Range("A7").Value = "Test accuracy on task 2 after task 2"
$ws.Range("B7").Value = 88.06
$ws.Range("C7").Value = 91.5
$ws.Range("D7").Value = 87.93000000000001
$ws.Range("E7").Value = 87.42
$ws.Range("F7").Value = 86.3
$ws.Range("G7").Value = 72.55
$ws.Range("H7").Value = 80.06999999999999

# Row 8: Test average accuracy after task 2
$ws.Range("A8").Value = "Test average accuracy after task 2"
$ws.Range("B8").Value = 49.155
$ws.Range("C8").Value = 91.5
$ws.Range("D8").Value = 92.16
$ws.Range("E8").Value = 92.59999999999999
$ws.Range("F8").Value = 92.205
$ws.Range("G8").Value = 63.915
$ws.Range("H8").Value = 62.27499999999999
